$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the account summary totals (base de datos EC actualizada) ---
# VALOR MORA total
$ws.Range("E11").Value = 113880
# Cant. Periodos
$ws.Range("F13").Value = 2

# --- Add a new "estado de cuenta" detail row (parte 1 de nuevos estado de cuenta) ---
# Duplicate the existing worker/period row (row 16) as a new row 17, keeping the
# same look & feel (font, fills, borders, number formats) as the existing data row.
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert()

# Re-apply the number/border formatting of the trailing blank cells (H:J) which the
# insert otherwise leaves unstyled.
$ws.Range("H16:J16").Copy()
$ws.Range("H17:J17").PasteSpecial(-4122)

# New row is for period 2508 (same worker, new mora period)
$ws.Range("E17").Value = "2508"

$excel.CutCopyMode = 0
